$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 27: new ESTIM log entry (RS neuron, 2018.08.21, 16:54:25)
$ws.Cells.Item(27, 1).Value = "'2018.08.21"
$ws.Cells.Item(27, 1).Style = "Normal"
$ws.Cells.Item(27, 2).Value = "16:54:25"
$ws.Cells.Item(27, 3).Value = "RS"
$ws.Cells.Item(27, 4).Value = 10
$ws.Cells.Item(27, 5).Value = 250
$ws.Cells.Item(27, 6).Value = 0.1
$ws.Cells.Item(27, 7).Value = 0.97
$ws.Cells.Item(27, 8).Value = 3475
$ws.Cells.Item(27, 9).Value = 0.4
$ws.Cells.Item(27, 10).Value = 0
$ws.Cells.Item(27, 11).Value = "N/A"
$ws.Cells.Item(27, 12).Value = "N/A"

# Row 28: new ESTIM log entry (RS neuron, 2018.08.21, 16:56:36)
$ws.Cells.Item(28, 1).Value = "'2018.08.21"
$ws.Cells.Item(28, 1).Style = "Normal"
$ws.Cells.Item(28, 2).Value = "16:56:36"
$ws.Cells.Item(28, 3).Value = "RS"
$ws.Cells.Item(28, 4).Value = 10
$ws.Cells.Item(28, 5).Value = 250
$ws.Cells.Item(28, 6).Value = 0.1
$ws.Cells.Item(28, 7).Value = 0.97
$ws.Cells.Item(28, 8).Value = 3475
$ws.Cells.Item(28, 9).Value = 0.41
$ws.Cells.Item(28, 10).Value = 0
$ws.Cells.Item(28, 11).Value = "N/A"
$ws.Cells.Item(28, 12).Value = "N/A"
